# Atualização automática de preços de eletricidade
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Day (serial date number, formatted as YYYY-MM-DD via existing cell style)
$ws.Range("A2").Value = 46044

# Hourly prices 0h-1h .. 23h-24h (columns B..Y)
$ws.Range("B2").Value = 41.06
$ws.Range("C2").Value = 33.72
$ws.Range("D2").Value = 44.4
$ws.Range("E2").Value = 34.37
$ws.Range("F2").Value = 30.64
$ws.Range("G2").Value = 43.43
$ws.Range("H2").Value = 48.14
$ws.Range("I2").Value = 56.23
$ws.Range("J2").Value = 75.59999999999999
$ws.Range("K2").Value = 82.31999999999999
$ws.Range("L2").Value = 74.06
$ws.Range("M2").Value = 56.7
$ws.Range("N2").Value = 53.31
$ws.Range("O2").Value = 53.34
$ws.Range("P2").Value = 55.06
$ws.Range("Q2").Value = 63.99
$ws.Range("R2").Value = 68.8
$ws.Range("S2").Value = 78.55
$ws.Range("T2").Value = 80.73999999999999
$ws.Range("U2").Value = 81.83
$ws.Range("V2").Value = 80.98
$ws.Range("W2").Value = 77.93000000000001
$ws.Range("X2").Value = 65.83
$ws.Range("Y2").Value = 55.06

# Price_Daily_Avg
$ws.Range("Z2").Value = 59.84

# Slot_4h_max / Slot_4h_price
$ws.Range("AA2").Value = "16h-20h"
$ws.Range("AB2").Value = 77.48

# Slot_2h_frist / Slot_2h_frist_price
$ws.Range("AC2").Value = "18h-20h"
$ws.Range("AD2").Value = 81.28

# Slot_2h_second / Slot_2h_second_price
$ws.Range("AE2").Value = "20h-22h"
$ws.Range("AF2").Value = 79.45999999999999

# Slot_min_price
$ws.Range("AG2").Value = "0h-23h"
